# The author renamed the first worksheet (fixing its casing) and switched
# the workbook so that sheet is the active/selected one on open, instead of
# the second sheet.

$wb = $excel.ActiveWorkbook

# Rename "Input_config" -> "input_config"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "input_config"

# Make "input_config" the active sheet (it was "input_config_2" before)
$ws1.Activate()
